# Insert a new price record at row 269 (weekly update), pushing the
# existing rows 269-331 down to 270-332.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 269, shifting the rest of the table down.
$ws.Rows.Item(269).Insert()

# Populate the new row with the latest week's observation.
$ws.Cells.Item(269, 1).Value = 5
$ws.Cells.Item(269, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(269, 3).Value = "Maule"
$ws.Cells.Item(269, 4).Value = 44943
$ws.Cells.Item(269, 5).Value = 7
$ws.Cells.Item(269, 6).Value = "Fruta"
$ws.Cells.Item(269, 7).Value = 100108
$ws.Cells.Item(269, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(269, 9).Value = 100108005
$ws.Cells.Item(269, 10).Value = "Piña"
$ws.Cells.Item(269, 11).Value = "Caramelo"
$ws.Cells.Item(269, 12).Value = "Segunda"
$ws.Cells.Item(269, 13).Value = 250
$ws.Cells.Item(269, 14).Value = 17000
$ws.Cells.Item(269, 15).Value = 17000
$ws.Cells.Item(269, 16).Value = 17000
$ws.Cells.Item(269, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(269, 18).Value = "Ecuador"
$ws.Cells.Item(269, 19).Value = 1214
$ws.Cells.Item(269, 20).Value = 14
